$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.AutoFit() | Out-Null
